# Actualización desde MV -datos-
# Adds the next 7 trading days of "Compra de dólares por licitación" data
# (08-09-2021 .. 16-09-2021) as new rows 164-170 on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Columns: A=Serie (date, stored as text), B=Monto licitado,
# C=Oferta, D=Monto adjudicado, E=Tipo cambio promedio ponderado adjudicado
$data = @(
    @("08-09-2021", 40, 175, 40, 788),
    @("09-09-2021", 40, 149, 40, 791),
    @("10-09-2021", 40, 102, 40, 790),
    @("13-09-2021", 40, 136, 40, 785),
    @("14-09-2021", 40, 145, 40, 782),
    @("15-09-2021", 40, 149, 40, 780),
    @("16-09-2021", 40, 101, 40, 782)
)

$startRow = 164
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $vals = $data[$i]

    # The "Serie" dates are day-first (dd-mm-yyyy) text, same as every
    # existing row above. Some of them (08/09/10-09-2021) are ambiguous
    # and would otherwise be auto-recognized as mm-dd-yyyy dates, so enter
    # them with a leading apostrophe to force plain text, then strip the
    # resulting quote-prefix formatting so the cell keeps the workbook's
    # default (unstyled) look, matching the rest of the column.
    $cA = $ws.Cells.Item($row, 1)
    $cA.Value = "'" + $vals[0]
    $cA.ClearFormats()

    $ws.Cells.Item($row, 2).Value = $vals[1]
    $ws.Cells.Item($row, 3).Value = $vals[2]
    $ws.Cells.Item($row, 4).Value = $vals[3]
    $ws.Cells.Item($row, 5).Value = $vals[4]
}
